# Auto-generated Excel COM-interop script applying the scheduled market-data refresh.
# Updates cached price/profit figures (columns H:N) across multiple sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2226.577
$ws.Range("I129").Value = 1123.6923
$ws.Range("J129").Value = 3329.4614
$ws.Range("K129").Value = 3371.0769
$ws.Range("L129").Value = 9988.3842
$ws.Range("M129").Value = 1628.9231
$ws.Range("N129").Value = -19988.3842
$ws.Range("H135").Value = 804.5517
$ws.Range("I135").Value = 567.85187
$ws.Range("K135").Value = 5110.66683
$ws.Range("M135").Value = -2575.66683
$ws.Range("H136").Value = 53390
$ws.Range("J136").Value = 53390
$ws.Range("L136").Value = 53390
$ws.Range("N136").Value = -63590

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 84784.8
$ws.Range("I28").Value = 66641.336
$ws.Range("J28").Value = 112000
$ws.Range("K28").Value = 66641.336
$ws.Range("L28").Value = 112000
$ws.Range("M28").Value = -66449.336
$ws.Range("N28").Value = -112384
$ws.Range("H88").Value = 980.2222
$ws.Range("I88").Value = 865
$ws.Range("J88").Value = 1037.8334
$ws.Range("K88").Value = 865
$ws.Range("L88").Value = 1037.8334
$ws.Range("M88").Value = -459
$ws.Range("N88").Value = -1849.8334
$ws.Range("H91").Value = 980.2222
$ws.Range("I91").Value = 865
$ws.Range("J91").Value = 1037.8334
$ws.Range("K91").Value = 865
$ws.Range("L91").Value = 1037.8334
$ws.Range("M91").Value = 539
$ws.Range("N91").Value = -3845.8334
$ws.Range("H97").Value = 1145.5294
$ws.Range("I97").Value = 764.93335
$ws.Range("K97").Value = 764.93335
$ws.Range("M97").Value = -268.93335
$ws.Range("H99").Value = 84784.8
$ws.Range("I99").Value = 66641.336
$ws.Range("J99").Value = 112000
$ws.Range("K99").Value = 66641.336
$ws.Range("L99").Value = 112000
$ws.Range("M99").Value = -63646.336
$ws.Range("N99").Value = -117990
$ws.Range("H122").Value = 7052.5264
$ws.Range("I122").Value = 6153.6924
$ws.Range("K122").Value = 18461.0772
$ws.Range("M122").Value = -16011.0772
$ws.Range("H128").Value = 99944.5
$ws.Range("J128").Value = 99944.5
$ws.Range("L128").Value = 99944.5
$ws.Range("N128").Value = -109904.5
$ws.Range("H133").Value = 95027.75
$ws.Range("J133").Value = 89666.664
$ws.Range("L133").Value = 89666.664
$ws.Range("N133").Value = -94726.664
$ws.Range("H138").Value = 103331.336
$ws.Range("J138").Value = 109995
$ws.Range("L138").Value = 109995
$ws.Range("N138").Value = -120275

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 101785.71
$ws.Range("I26").Value = 100000
$ws.Range("J26").Value = 102083.336
$ws.Range("K26").Value = 100000
$ws.Range("L26").Value = 102083.336
$ws.Range("M26").Value = -99708
$ws.Range("N26").Value = -102667.336
$ws.Range("H28").Value = 58245.5
$ws.Range("J28").Value = 58245.5
$ws.Range("L28").Value = 58245.5
$ws.Range("N28").Value = -58833.5
$ws.Range("H40").Value = 39529
$ws.Range("J40").Value = 39529
$ws.Range("L40").Value = 39529
$ws.Range("N40").Value = -40059
$ws.Range("H86").Value = 3635.3103
$ws.Range("I86").Value = 3407.3684
$ws.Range("J86").Value = 4068.4
$ws.Range("K86").Value = 3407.3684
$ws.Range("L86").Value = 4068.4
$ws.Range("M86").Value = -2284.3684
$ws.Range("N86").Value = -6314.4
$ws.Range("H89").Value = 3635.3103
$ws.Range("I89").Value = 3407.3684
$ws.Range("J89").Value = 4068.4
$ws.Range("K89").Value = 17036.842
$ws.Range("L89").Value = 20342
$ws.Range("M89").Value = -11420.842
$ws.Range("N89").Value = -31574
$ws.Range("H98").Value = 111500
$ws.Range("J98").Value = 111500
$ws.Range("L98").Value = 111500
$ws.Range("N98").Value = -117490
$ws.Range("H108").Value = 109499.5
$ws.Range("J108").Value = 109499.5
$ws.Range("L108").Value = 109499.5
$ws.Range("N108").Value = -117179.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 239.66667
$ws.Range("I7").Value = 120
$ws.Range("J7").Value = 389.25
$ws.Range("K7").Value = 120
$ws.Range("L7").Value = 389.25
$ws.Range("M7").Value = -7
$ws.Range("N7").Value = -615.25
$ws.Range("H75").Value = 96498.5
$ws.Range("J75").Value = 102798.6
$ws.Range("L75").Value = 102798.6
$ws.Range("N75").Value = -104794.6
$ws.Range("H78").Value = 96498.5
$ws.Range("J78").Value = 102798.6
$ws.Range("L78").Value = 308395.8
$ws.Range("N78").Value = -318379.8
$ws.Range("H86").Value = 3444
$ws.Range("I86").Value = 3533.25
$ws.Range("K86").Value = 3533.25
$ws.Range("M86").Value = -2410.25
$ws.Range("H89").Value = 3444
$ws.Range("I89").Value = 3533.25
$ws.Range("K89").Value = 17666.25
$ws.Range("M89").Value = -12050.25
$ws.Range("H122").Value = 2368.2
$ws.Range("I122").Value = 2368.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7104.599999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4654.599999999999
$ws.Range("N122").ClearContents()
$ws.Range("H141").Value = 695132.7
$ws.Range("I141").Value = 130000
$ws.Range("J141").Value = 757925.25
$ws.Range("K141").Value = 130000
$ws.Range("L141").Value = 757925.25
$ws.Range("M141").Value = -124820
$ws.Range("N141").Value = -768285.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1410.5555
$ws.Range("I102").Value = 1297.3334
$ws.Range("K102").Value = 1297.3334
$ws.Range("M102").Value = 324.6666
$ws.Range("H126").Value = 1157
$ws.Range("I126").Value = 1190
$ws.Range("K126").Value = 3570
$ws.Range("M126").Value = -1100
$ws.Range("H128").Value = 152780
$ws.Range("J128").Value = 152780
$ws.Range("L128").Value = 152780
$ws.Range("N128").Value = -162740
$ws.Range("H132").Value = 1988.9722
$ws.Range("I132").Value = 1685.8462
$ws.Range("K132").Value = 5057.5386
$ws.Range("M132").Value = -2527.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3849.75
$ws.Range("J68").Value = 3849.75
$ws.Range("L68").Value = 3849.75
$ws.Range("N68").Value = -5347.75
$ws.Range("H71").Value = 3849.75
$ws.Range("J71").Value = 3849.75
$ws.Range("L71").Value = 19248.75
$ws.Range("N71").Value = -26736.75
$ws.Range("H93").Value = 1180.2941
$ws.Range("I93").Value = 990.9167
$ws.Range("K93").Value = 990.9167
$ws.Range("M93").Value = 257.0833
$ws.Range("H122").Value = 7121.1113
$ws.Range("J122").Value = 4147.5
$ws.Range("L122").Value = 12442.5
$ws.Range("N122").Value = -17342.5
$ws.Range("H132").Value = 3408.4243
$ws.Range("I132").Value = 2439.2778
$ws.Range("K132").Value = 7317.8334
$ws.Range("M132").Value = -4787.8334
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("N134").Value = -110140
$ws.Range("H139").Value = 97701
$ws.Range("J139").Value = 97701
$ws.Range("L139").Value = 97701
$ws.Range("N139").Value = -107981

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 189439.75
$ws.Range("J128").Value = 189439.75
$ws.Range("L128").Value = 189439.75
$ws.Range("N128").Value = -199399.75
$ws.Range("H132").Value = 4662.7085
$ws.Range("I132").Value = 3595.8823
$ws.Range("K132").Value = 10787.6469
$ws.Range("M132").Value = -8257.6469

Write-Host "Applied 191 cell updates and 2 clears across 7 sheets."
